$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing subject codes (Cod Materia column, A)
# Write order is chosen to reproduce the exact shared-string append order
$ws.Range("A4").Value = "ACO"
$ws.Range("A3").Value = "FI1"
$ws.Range("A6").Value = "SYO"
$ws.Range("A5").Value = "MAD"

# Add two new rows of class schedule data, copying the row-7 formatting first
# so the new rows share the same styles (no new style entries created)
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A9:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "ADS"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Lunes"
$ws.Range("D8").Value = 0.45833333333333331
$ws.Range("E8").Value = 25

$ws.Range("A9").Value = "IYS"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Lunes"
$ws.Range("D9").Value = 0.45833333333333331
$ws.Range("E9").Value = 25

# Update selection to mirror the saved cursor position
$ws.Range("H7").Select()
